$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column L: "Authorship Resource" ---
$ws.Range("L1").Value = "Authorship Resource"
$ws.Range("L2:L25").Value = "Daniela Subotic, Noémi Villars-Amberg"

# --- Column width adjustments ---
# Column G (7): 119.5 -> 16.83203125
$ws.Columns.Item(7).ColumnWidth = 16
# Column J (10): 78.5 -> 39.5
$ws.Columns.Item(10).ColumnWidth = 38.67
# Column K (11): 74.1640625 (bestFit) -> 35.6640625
$ws.Columns.Item(11).ColumnWidth = 34.83
# Column L (12, new): -> 92.5
$ws.Columns.Item(12).ColumnWidth = 91.67

# --- Sheet view / window state ---
$win = $excel.ActiveWindow
$win.Zoom = 149
$null = $ws.Range("L2:L25").Select()
